$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3: new cell that keeps the ORIGINAL look of A1 (Calibri, italic, 24pt) ---
# Capture A1's current formatting (before we touch it) and apply it to B3 as a
# single "paste formats" operation, so the existing italic/24pt font is reused
# rather than rebuilt property-by-property.
$ws.Range("B3").Value = "24 pt Italic"
$ws.Range("A1").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # -4122 = xlPasteFormats

# --- A1: new text, new look (Times New Roman, bold) ---
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = "Bold Times New Roman"
$ws.Range("A1").Font.Name = "Times New Roman"
$ws.Range("A1").Font.Bold = $true
